$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# 1. "Primitive source control..." -> "Basic source control..."
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Basic source control: copy the entire folder and name it " + [char]8220 + "OLD_v3" + [char]8221

# 2. New bullet after it (lvl 1): "Better than nothing, but quickly becomes difficult to work with"
$para1.InsertAfter("`rBetter than nothing, but quickly becomes difficult to work with")
$para2 = $tr.Paragraphs(2, 1)
$para2.IndentLevel = 2

# 3. Blank spacer paragraph (lvl 1, no bullet) right after it
$para2.InsertAfter("`r")
$para3 = $tr.Paragraphs(3, 1)
$para3.IndentLevel = 2
$para3.ParagraphFormat.Bullet.Visible = 0

# 4. Another blank spacer paragraph after "Good for tracking a mix of different file types"
$para5 = $tr.Paragraphs(5, 1)
$para5.InsertAfter("`r")
$para6 = $tr.Paragraphs(6, 1)
$para6.IndentLevel = 2
$para6.ParagraphFormat.Bullet.Visible = 0
